$wb = $excel.ActiveWorkbook

# ---------- Sheet1 ----------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 3: ConnectStatus flips from Yes to No
$ws1.Cells.Item(3,6).Value = "No"

# Row 4: CurrentPage + LastDate get populated, ConnectStatus flips No -> Yes
$ws1.Cells.Item(4,4).Value = 2
$ws1.Cells.Item(4,5).Value = 44831.957442129627
$ws1.Cells.Item(4,6).Value = "Yes"

# Row 5: CurrentPage + LastDate get populated, ConnectStatus flips No -> Yes
$ws1.Cells.Item(5,4).Value = 2
$ws1.Cells.Item(5,5).NumberFormat = "m/d/yy h:mm"
$ws1.Cells.Item(5,5).Value = 44831.959918981483
$ws1.Cells.Item(5,6).Value = "Yes"

# New rows 6-10: additional SearchKeyword entries for Prashant Mewada / Bengaluru
$ws1.Cells.Item(6,1).Value = "Prashant Mewada"
$ws1.Cells.Item(6,2).Value = "HR Executive"
$ws1.Cells.Item(6,3).Value = "Bengaluru"

$ws1.Cells.Item(7,1).Value = "Prashant Mewada"
$ws1.Cells.Item(7,2).Value = "Recruitment Specialist"
$ws1.Cells.Item(7,3).Value = "Bengaluru"

$ws1.Cells.Item(8,1).Value = "Prashant Mewada"
$ws1.Cells.Item(8,2).Value = "IT recruiter"
$ws1.Cells.Item(8,3).Value = "Bengaluru"

$ws1.Cells.Item(9,1).Value = "Prashant Mewada"
$ws1.Cells.Item(9,2).Value = "Information Technology recruiter"
$ws1.Cells.Item(9,3).Value = "Bengaluru"

$ws1.Cells.Item(10,1).Value = "Prashant Mewada"
$ws1.Cells.Item(10,2).Value = "HR recruiter"
$ws1.Cells.Item(10,3).Value = "Bengaluru"

# Column B widened to fit the longest SearchKeyword text
$ws1.Columns.Item(2).ColumnWidth = 30.1666666

# ---------- Sheet2 ----------
$ws2 = $wb.Worksheets.Item("Sheet2")

function Add-Sheet2Row {
    param($r, $name, $keyword, $location, $status, $url, $timestamp)
    $ws2.Cells.Item($r,1).Value = $name
    $ws2.Cells.Item($r,2).Value = $keyword
    $ws2.Cells.Item($r,3).Value = $location
    $ws2.Cells.Item($r,4).Value = $status
    $ws2.Cells.Item($r,5).Value = $url
    $ws2.Cells.Item($r,6).NumberFormat = "m/d/yy h:mm"
    $ws2.Cells.Item($r,6).Value = $timestamp
}

Add-Sheet2Row 12 "Lalit .." "technical recruiter" "Bengaluru" "Connect - Sent" "https://www.linkedin.com/in/lalit-4bba7924a?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAD29Do8B6n63DBAgEUbzASH-DNU18DFXrgE" 44831.956238425926
Add-Sheet2Row 13 "Kavya Shree" "technical recruiter" "Bengaluru" "Connect - Sent" "https://www.linkedin.com/in/kavyashreeav?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAAt7G7gBZBipGj3co0ddZiNKK__1xhLBVTQ" 44831.956284722219
Add-Sheet2Row 14 "Shwetha S L" "technical recruiter" "Bengaluru" "Pending - Already Sent" "https://www.linkedin.com/in/shwetha-s-l-7904a924?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAAUOWxwB8e7Wsf6VKyY9g0OR3UwFmfLzGHA" 44831.956284722219
Add-Sheet2Row 15 "Harshith K" "technical recruiter" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/harshith-k-455a3224b?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAD3pKwEBmHbKU2rk6XEXGbjoPGIUjyrjH8E" 44831.95653935185
Add-Sheet2Row 16 "Hemant S" "technical recruiter" "Bengaluru" "Connect - Sent" "https://www.linkedin.com/in/hemant-s-0a1a481b7?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAADKFonkBhhq0LQIoBZx7pU8v-JOFtIrAkGk" 44831.956585648149
Add-Sheet2Row 17 "R Aarthi She/Her" "technical recruiter" "Bengaluru" "Follow - Sent" "https://www.linkedin.com/in/aarthirecruitergoogle?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABmDMYcBHai2bR4gJ9Hy7ATnN6Mv9Jrb2mA" 44831.95684027778
Add-Sheet2Row 18 "Ankita Goel" "technical recruiter" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/ankita-goel-007ab5147?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACOUhMcB24kw_osHDEw3La9F4hBGBPNpUWI" 44831.957083333335
Add-Sheet2Row 19 "Srinath M" "technical recruiter" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/srinath-m-546651244?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAADybpsYBki2f6h7c-j4Dy7p2SoaKeUEjNWs" 44831.957326388889
Add-Sheet2Row 20 "Vidya hj" "technical recruiter" "Bengaluru" "Connect - Sent" "https://www.linkedin.com/in/vidyahj?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABsAHJkB7r-ikZXuMH8q9Opw2TYU_pDzVdE" 44831.957372685189
Add-Sheet2Row 21 "Anuraag Mandanna" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/anuraag-mandanna-a06b75158?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACXnYmYBc_wNWufnIVXRQLIb-_UwY9OuUCI" 44831.958101851851
Add-Sheet2Row 22 "Ramya D" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/ramya-d-b7310b106?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABq2fMcBApC11S5AGrSeSGKVEsHnOWMATAY" 44831.958333333336
Add-Sheet2Row 23 "Apeksha ." "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/apeksha3?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAAnEoRABgpnQUOiRjkIujbUbbLCDgTNMkLw" 44831.958599537036
Add-Sheet2Row 24 "Samreen Jabbar" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/samreen-jabbar-220239172?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACj4KPgB1-GRuUKp7Y1JW6EGeCZpfmZcxLU" 44831.95884259259
Add-Sheet2Row 25 "Astha Handa" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/asthahanda?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACYJDAkB6F5meb0X1a8RpFMJ736SaBMosIE" 44831.959097222221
Add-Sheet2Row 26 "Preethika gk" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/preethika-gk-b60144140?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACIxMXgBCHQlrQ2Kd6vWjCNNv_NV0HjrNr0" 44831.959340277775
Add-Sheet2Row 27 "Sananda Basu" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/sananda-basu-666a3358?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAAxBQu4B1pbluJLT08oDB_nv5r7FIefoIlg" 44831.95957175926
Add-Sheet2Row 28 "Luzanna Virgina Barretto" "Human Resource" "Bengaluru" "Connect - Sent" "https://www.linkedin.com/in/luzanna-virgina-barretto-b170401a2?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAC9Y6kMBfuv5q1qNF5jMAsBpaz55R6SS8XU" 44831.959618055553
Add-Sheet2Row 29 "Athiba Syed" "Human Resource" "Bengaluru" "Message - Sent" "https://www.linkedin.com/in/athiba-syed-432923176?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACnIhKYBbKNHfAdCeQfQs1p6FsmA0198ZYs" 44831.959872685184

# ---------- Selections ----------
$ws2.Range("G2").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A2:A9").Select()

$ws1.Activate()
$ws1.Range("F6").Select()
